# Update "Forecast Comparison" sheet: Prophet Forecast (B) and yhat_upper (D) columns
$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary = $wb.Worksheets.Item("Summary")

# Row => (B new value or $null if unchanged, D new value)
$updates = @{
    2  = @{ B = $null; D = 190.0910339471775 }
    3  = @{ B = 207;   D = 233.3275293256467 }
    4  = @{ B = 212;   D = 240.8982021110415 }
    5  = @{ B = 177;   D = 203.2885980182523 }
    6  = @{ B = 137;   D = 162.8886632926573 }
    7  = @{ B = $null; D = 146.8522370952435 }
    8  = @{ B = $null; D = 146.537376583008 }
    9  = @{ B = 116;   D = 141.0210774705989 }
    10 = @{ B = $null; D = 136.1640153028027 }
    11 = @{ B = 119;   D = 146.252532241961 }
    12 = @{ B = $null; D = 175.5692787002968 }
    13 = @{ B = $null; D = 200.029598337419 }
    14 = @{ B = 149;   D = 176.6062986051131 }
    15 = @{ B = 87;    D = 116.1634683976632 }
    16 = @{ B = $null; D = 68.45867424498672 }
    17 = @{ B = $null; D = 83.86851019117597 }
    18 = @{ B = $null; D = 144.5568439352778 }
    19 = @{ B = $null; D = 176.2469300397337 }
    20 = @{ B = $null; D = 152.7115243547536 }
    21 = @{ B = $null; D = 103.5775184397716 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.B) {
        $wsForecast.Range("B$row").Value = $vals.B
    }
    $wsForecast.Range("D$row").Value = $vals.D
}

# Update "Summary" sheet totals (stored as text strings, not numbers).
# Temporarily switch the cell to a text NumberFormat so the COM layer keeps
# the assigned value as a string instead of auto-coercing it to a number,
# then restore the cell's original style so no visible formatting changes.
function Set-TextValue($range, [string]$text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $wsSummary.Range("B9") "2143"
Set-TextValue $wsSummary.Range("B10") "1253"
Set-TextValue $wsSummary.Range("B11") "760"
Set-TextValue $wsSummary.Range("B12") "212"
